# Update scripts with new TPM (transcripts per million) derived NATMI
# ligand/receptor values for the Hbegf-Cd9 pair. Re-scored ligand/receptor
# expression (G,H,M,N) and their downstream derived-specificity and
# edge-weight columns (I,J,O,P,Q,R,S,T) for rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.723979
$ws.Range("H2").Value = 23.171937
$ws.Range("I2").Value = 0.471042132528101
$ws.Range("J2").Value = 0.471042132528101
$ws.Range("M2").Value = 19.515399
$ws.Range("N2").Value = 58.546197
$ws.Range("O2").Value = 0.1046357846766865
$ws.Range("P2").Value = 0.1046357846766865
$ws.Range("Q2").Value = 150.736532052621
$ws.Range("R2").Value = 1356.628788473589
$ws.Range("S2").Value = 0.04928786315285762
$ws.Range("T2").Value = 0.04928786315285762
$ws.Range("G3").Value = 7.723979
$ws.Range("H3").Value = 23.171937
$ws.Range("I3").Value = 0.471042132528101
$ws.Range("J3").Value = 0.471042132528101
$ws.Range("M3").Value = 133.0753813333333
$ws.Range("N3").Value = 399.226144
$ws.Range("O3").Value = 0.7135107484588257
$ws.Range("P3").Value = 0.7135107484588257
$ws.Range("Q3").Value = 1027.871450835659
$ws.Range("R3").Value = 9250.843057520928
$ws.Range("S3").Value = 0.3360936245357667
$ws.Range("T3").Value = 0.3360936245357667
$ws.Range("G4").Value = 7.723979
$ws.Range("H4").Value = 23.171937
$ws.Range("I4").Value = 0.471042132528101
$ws.Range("J4").Value = 0.471042132528101
$ws.Range("O4").Value = 0.1818534668644878
$ws.Range("P4").Value = 0.1818534668644878
$ws.Range("Q4").Value = 261.9750119101143
$ws.Range("R4").Value = 2357.775107191029
$ws.Range("S4").Value = 0.08566064483947668
$ws.Range("T4").Value = 0.08566064483947668
$ws.Range("I5").Value = 0.2460132574367717
$ws.Range("J5").Value = 0.2460132574367717
$ws.Range("M5").Value = 19.515399
$ws.Range("N5").Value = 58.546197
$ws.Range("O5").Value = 0.1046357846766865
$ws.Range("P5").Value = 0.1046357846766865
$ws.Range("Q5").Value = 78.72583513062999
$ws.Range("R5").Value = 708.53251617567
$ws.Range("S5").Value = 0.0257417902327643
$ws.Range("T5").Value = 0.0257417902327643
$ws.Range("I6").Value = 0.2460132574367717
$ws.Range("J6").Value = 0.2460132574367717
$ws.Range("M6").Value = 133.0753813333333
$ws.Range("N6").Value = 399.226144
$ws.Range("O6").Value = 0.7135107484588257
$ws.Range("P6").Value = 0.7135107484588257
$ws.Range("Q6").Value = 536.8309677293155
$ws.Range("R6").Value = 4831.478709563839
$ws.Range("S6").Value = 0.1755331034445047
$ws.Range("T6").Value = 0.1755331034445047
$ws.Range("I7").Value = 0.2460132574367717
$ws.Range("J7").Value = 0.2460132574367717
$ws.Range("O7").Value = 0.1818534668644878
$ws.Range("P7").Value = 0.1818534668644878
$ws.Range("S7").Value = 0.04473836375950267
$ws.Range("T7").Value = 0.04473836375950267
$ws.Range("I8").Value = 0.2829446100351274
$ws.Range("J8").Value = 0.2829446100351274
$ws.Range("M8").Value = 19.515399
$ws.Range("N8").Value = 58.546197
$ws.Range("O8").Value = 0.1046357846766865
$ws.Range("P8").Value = 0.1046357846766865
$ws.Range("Q8").Value = 90.54410706484299
$ws.Range("R8").Value = 814.896963583587
$ws.Range("S8").Value = 0.02960613129106464
$ws.Range("T8").Value = 0.02960613129106464
$ws.Range("I9").Value = 0.2829446100351274
$ws.Range("J9").Value = 0.2829446100351274
$ws.Range("M9").Value = 133.0753813333333
$ws.Range("N9").Value = 399.226144
$ws.Range("O9").Value = 0.7135107484588257
$ws.Range("P9").Value = 0.7135107484588257
$ws.Range("Q9").Value = 617.4196886848248
$ws.Range("R9").Value = 5556.777198163423
$ws.Range("S9").Value = 0.2018840204785543
$ws.Range("T9").Value = 0.2018840204785543
$ws.Range("I10").Value = 0.2829446100351274
$ws.Range("J10").Value = 0.2829446100351274
$ws.Range("O10").Value = 0.1818534668644878
$ws.Range("P10").Value = 0.1818534668644878
$ws.Range("S10").Value = 0.05145445826550847
$ws.Range("T10").Value = 0.05145445826550847
